$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 554.34784
$ws.Cells.Item(17, 9).Value = 93.333336
$ws.Cells.Item(17, 10).Value = 569.8876299999999
$ws.Cells.Item(17, 11).Value = 280.000008
$ws.Cells.Item(17, 12).Value = 1709.66289
$ws.Cells.Item(17, 13).Value = -112.000008
$ws.Cells.Item(17, 14).Value = -2045.66289
$ws.Cells.Item(98, 8).Value = 1877.804
$ws.Cells.Item(98, 9).Value = 2282.925
$ws.Cells.Item(98, 10).Value = 404.63635
$ws.Cells.Item(98, 11).Value = 2282.925
$ws.Cells.Item(98, 12).Value = 404.63635
$ws.Cells.Item(98, 13).Value = -784.9250000000002
$ws.Cells.Item(98, 14).Value = -3400.63635
$ws.Cells.Item(107, 8).Value = 1282.4706
$ws.Cells.Item(107, 9).Value = 1405.7778
$ws.Cells.Item(107, 10).Value = 1143.75
$ws.Cells.Item(107, 11).Value = 1405.7778
$ws.Cells.Item(107, 12).Value = 1143.75
$ws.Cells.Item(107, 13).Value = 514.2221999999999
$ws.Cells.Item(107, 14).Value = -4983.75
$ws.Cells.Item(113, 8).Value = 1543.2667
$ws.Cells.Item(113, 9).Value = 1337.25
$ws.Cells.Item(113, 10).Value = 1618.1818
$ws.Cells.Item(113, 11).Value = 1337.25
$ws.Cells.Item(113, 12).Value = 1618.1818
$ws.Cells.Item(113, 13).Value = 1916.75
$ws.Cells.Item(113, 14).Value = -8126.1818
$ws.Cells.Item(122, 8).Value = 1877.804
$ws.Cells.Item(122, 9).Value = 2282.925
$ws.Cells.Item(122, 10).Value = 404.63635
$ws.Cells.Item(122, 11).Value = 6848.775000000001
$ws.Cells.Item(122, 12).Value = 1213.90905
$ws.Cells.Item(122, 13).Value = -4398.775000000001
$ws.Cells.Item(122, 14).Value = -6113.90905
$ws.Cells.Item(132, 8).Value = 280847.66
$ws.Cells.Item(132, 9).Value = 336827.12
$ws.Cells.Item(132, 10).Value = 950.3333
$ws.Cells.Item(132, 11).Value = 1010481.36
$ws.Cells.Item(132, 12).Value = 2850.9999
$ws.Cells.Item(132, 13).Value = -1007951.36
$ws.Cells.Item(132, 14).Value = -7910.9999
$ws.Cells.Item(135, 8).Value = 1440.8864
$ws.Cells.Item(135, 9).Value = 326.87097
$ws.Cells.Item(135, 10).Value = 4097.385
$ws.Cells.Item(135, 11).Value = 2941.83873
$ws.Cells.Item(135, 12).Value = 36876.465
$ws.Cells.Item(135, 13).Value = -406.8387299999999
$ws.Cells.Item(135, 14).Value = -41946.465
$ws.Cells.Item(141, 8).Value = 1104.0454
$ws.Cells.Item(141, 9).Value = 945.26666
$ws.Cells.Item(141, 11).Value = 2835.79998
$ws.Cells.Item(141, 13).Value = 2344.20002

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 9175.758
$ws.Cells.Item(32, 9).Value = 10085.186
$ws.Cells.Item(32, 11).Value = 10085.186
$ws.Cells.Item(32, 13).Value = -9798.186
$ws.Cells.Item(61, 8).Value = 20835754
$ws.Cells.Item(61, 9).Value = 22224670
$ws.Cells.Item(61, 10).Value = 2000
$ws.Cells.Item(61, 11).Value = 22224670
$ws.Cells.Item(61, 12).Value = 2000
$ws.Cells.Item(61, 13).Value = -22224458
$ws.Cells.Item(61, 14).Value = -2424
$ws.Cells.Item(74, 8).Value = 9098151
$ws.Cells.Item(74, 9).Value = 20000992
$ws.Cells.Item(74, 10).Value = 12451.25
$ws.Cells.Item(74, 11).Value = 20000992
$ws.Cells.Item(74, 12).Value = 12451.25
$ws.Cells.Item(74, 13).Value = -20000118
$ws.Cells.Item(74, 14).Value = -14199.25
$ws.Cells.Item(77, 8).Value = 9098151
$ws.Cells.Item(77, 9).Value = 20000992
$ws.Cells.Item(77, 10).Value = 12451.25
$ws.Cells.Item(77, 11).Value = 100004960
$ws.Cells.Item(77, 12).Value = 62256.25
$ws.Cells.Item(77, 13).Value = -100000592
$ws.Cells.Item(77, 14).Value = -70992.25
$ws.Cells.Item(101, 8).Value = 37500
$ws.Cells.Item(101, 10).Value = 37500
$ws.Cells.Item(101, 12).Value = 37500
$ws.Cells.Item(101, 14).Value = -43990
$ws.Cells.Item(122, 8).Value = 1365.6666
$ws.Cells.Item(122, 9).Value = 1262.2727
$ws.Cells.Item(122, 10).Value = 1650
$ws.Cells.Item(122, 11).Value = 3786.8181
$ws.Cells.Item(122, 12).Value = 4950
$ws.Cells.Item(122, 13).Value = -1336.8181
$ws.Cells.Item(122, 14).Value = -9850
$ws.Cells.Item(136, 8).Value = 20835754
$ws.Cells.Item(136, 9).Value = 22224670
$ws.Cells.Item(136, 10).Value = 2000
$ws.Cells.Item(136, 11).Value = 66674010
$ws.Cells.Item(136, 12).Value = 6000
$ws.Cells.Item(136, 13).Value = -66671460
$ws.Cells.Item(136, 14).Value = -11100

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(31, 8).Value = 1500
$ws.Cells.Item(31, 9).Value = 1500
$ws.Cells.Item(31, 11).Value = 1500
$ws.Cells.Item(31, 13).Value = -1248
$ws.Cells.Item(64, 8).Value = 1480.0952
$ws.Cells.Item(64, 9).Value = 1211.3334
$ws.Cells.Item(64, 10).Value = 1587.6
$ws.Cells.Item(64, 11).Value = 1211.3334
$ws.Cells.Item(64, 12).Value = 1587.6
$ws.Cells.Item(64, 13).Value = -986.3334
$ws.Cells.Item(64, 14).Value = -2037.6
$ws.Cells.Item(67, 8).Value = 1480.0952
$ws.Cells.Item(67, 9).Value = 1211.3334
$ws.Cells.Item(67, 10).Value = 1587.6
$ws.Cells.Item(67, 11).Value = 1211.3334
$ws.Cells.Item(67, 12).Value = 1587.6
$ws.Cells.Item(67, 13).Value = -431.3334
$ws.Cells.Item(67, 14).Value = -3147.6
$ws.Cells.Item(134, 8).Value = 22710.34
$ws.Cells.Item(134, 9).Value = 29276.223
$ws.Cells.Item(134, 10).Value = 1222
$ws.Cells.Item(134, 11).Value = 87828.66900000001
$ws.Cells.Item(134, 12).Value = 3666
$ws.Cells.Item(134, 13).Value = -85293.66900000001
$ws.Cells.Item(134, 14).Value = -8736

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(74, 8).Value = 3000
$ws.Cells.Item(74, 9).Value = 0
$ws.Cells.Item(74, 11).Value = 0
$ws.Cells.Item(74, 13).ClearContents()
$ws.Cells.Item(77, 8).Value = 3000
$ws.Cells.Item(77, 9).Value = 0
$ws.Cells.Item(77, 11).Value = 0
$ws.Cells.Item(77, 13).ClearContents()
$ws.Cells.Item(94, 8).Value = 3300
$ws.Cells.Item(94, 10).Value = 3300
$ws.Cells.Item(94, 12).Value = 9900
$ws.Cells.Item(94, 14).Value = -11252
$ws.Cells.Item(131, 8).Value = 3072.3555
$ws.Cells.Item(131, 10).Value = 2375.6943
$ws.Cells.Item(131, 12).Value = 7127.0829
$ws.Cells.Item(131, 14).Value = -17207.0829
$ws.Cells.Item(133, 8).Value = 7976.6665
$ws.Cells.Item(133, 9).Value = 5515
$ws.Cells.Item(133, 10).Value = 8680
$ws.Cells.Item(133, 11).Value = 16545
$ws.Cells.Item(133, 12).Value = 26040
$ws.Cells.Item(133, 13).Value = -11485
$ws.Cells.Item(133, 14).Value = -36160

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(19, 8).Value = 2860.8696
$ws.Cells.Item(19, 9).Value = 1266.6666
$ws.Cells.Item(19, 10).Value = 3100
$ws.Cells.Item(19, 11).Value = 1266.6666
$ws.Cells.Item(19, 12).Value = 3100
$ws.Cells.Item(19, 13).Value = -978.6666
$ws.Cells.Item(19, 14).Value = -3676

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(19, 8).Value = 1190
$ws.Cells.Item(19, 9).Value = 380
$ws.Cells.Item(19, 10).Value = 2000
$ws.Cells.Item(19, 11).Value = 380
$ws.Cells.Item(19, 12).Value = 2000
$ws.Cells.Item(19, 13).Value = -210
$ws.Cells.Item(19, 14).Value = -2340
$ws.Cells.Item(55, 8).Value = 194.53334
$ws.Cells.Item(55, 9).Value = 170.61539
$ws.Cells.Item(55, 10).Value = 350
$ws.Cells.Item(55, 11).Value = 170.61539
$ws.Cells.Item(55, 12).Value = 350
$ws.Cells.Item(55, 13).Value = 2.384610000000009
$ws.Cells.Item(55, 14).Value = -696
$ws.Cells.Item(132, 8).Value = 4778
$ws.Cells.Item(132, 9).Value = 5783.913
$ws.Cells.Item(132, 10).Value = 3417.0588
$ws.Cells.Item(132, 11).Value = 17351.739
$ws.Cells.Item(132, 12).Value = 10251.1764
$ws.Cells.Item(132, 13).Value = -14821.739
$ws.Cells.Item(132, 14).Value = -15311.1764

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(13, 8).Value = 0
$ws.Cells.Item(13, 9).Value = 0
$ws.Cells.Item(13, 10).Value = 0
$ws.Cells.Item(13, 11).Value = 0
$ws.Cells.Item(13, 12).Value = 0
$ws.Cells.Item(13, 13).ClearContents()
$ws.Cells.Item(13, 14).ClearContents()
$ws.Cells.Item(126, 8).Value = 2483.4707
$ws.Cells.Item(126, 9).Value = 3350
$ws.Cells.Item(126, 10).Value = 894.8333
$ws.Cells.Item(126, 11).Value = 10050
$ws.Cells.Item(126, 12).Value = 2684.4999
$ws.Cells.Item(126, 13).Value = -7580
$ws.Cells.Item(126, 14).Value = -7624.4999
$ws.Cells.Item(136, 8).Value = 12409.444
$ws.Cells.Item(136, 10).Value = 1156.4445
$ws.Cells.Item(136, 12).Value = 3469.3335
$ws.Cells.Item(136, 14).Value = -8569.333500000001

Write-Host "Edit complete"